# Generate Report for handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) for the second data row
# (the e900ff5d... file) on both the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-17 07:35:25"
$wsZhCn.Range("G3").Value = "2016-01-17 07:36:09"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-17 07:35:35"
$wsDeDe.Range("G3").Value = "2016-01-17 07:36:26"
